$d = $word.ActiveDocument

# The paragraph "<id>p052r_1</id>" is currently split across three runs:
#   run 1: "<id>"    (Courier New, color 7f6000, sz 18)
#   run 2: "p052r_1" (default formatting)
#   run 3: "</id>"   (Courier New, color 7f6000, sz 18)
# Collapse them into a single run carrying the combined text
# "<id>p052r_1</id>" (adopting run 1's formatting), matching the commit's
# "tcn" cleanup that merges the newly-downloaded id fragments.
$range = $d.Content
$range.Find.Execute("<id>p052r_1</id>", $true, $false, $false, $false, $false, `
                     $true, 1, $false, "<id>p052r_1</id>", 2)
